# Color LED 회로 수정 / 저항 추가
#   - footer "date" field text 2015-05-29 -> 2015-08-10 (master + every layout)
#   - slide 17: merge the two "Circuits for " / "Appendix" runs into one run
#   - slide 7 (Color LED Circuit): add three 330 ohm resistor picture+label
#     pairs, copied from the existing resistor on slide 3 (LED Circuit)

$p = $ppt.ActivePresentation
$EMU = 12700.0

# ---------------------------------------------------------------------
# 1) Footer date placeholder: 2015-05-29 -> 2015-08-10
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "2015-05-29") {
                # force an actual change so the engine re-serialises the run
                $sh.TextFrame.TextRange.Text = "2015-08-10x"
                $sh.TextFrame.TextRange.Text = "2015-08-10"
            }
        }
    }
}

Update-DateShapes $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 17 ("Circuits for " + "Appendix" -> single run)
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
for ($i = 1; $i -le $s17.Shapes.Count; $i++) {
    $sh = $s17.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "Circuits for Appendix") {
            $sh.TextFrame.TextRange.Text = "x"
            $sh.TextFrame.TextRange.Text = "Circuits for Appendix"
        }
    }
}

# ---------------------------------------------------------------------
# 3) Slide 7 ("Color LED Circuit"): add three 330 ohm resistors
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s7 = $p.Slides.Item(7)

$picSrc = $null
$txtSrc = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.Name -eq "Picture 5") { $picSrc = $sh }
    if ($sh.Name -eq "TextBox 24") { $txtSrc = $sh }
}

$sets = @(
    @{ picX = 2555776; picY = 2522375; txX = 2445032; txY = 2317696; txName = "TextBox 25" },
    @{ picX = 2599958; picY = 4156337; txX = 2489214; txY = 3951658; txName = "TextBox 29" },
    @{ picX = 2597432; picY = 4941316; txX = 2486688; txY = 4736637; txName = "TextBox 31" }
)

foreach ($set in $sets) {
    $picSrc.Copy()
    $pic = $s7.Shapes.Paste().Item(1)
    $pic.Name = "Picture 5"
    $pic.Left = $set.picX / $EMU
    $pic.Top = $set.picY / $EMU
    $pic.Width = 725690 / $EMU
    $pic.Height = 337677 / $EMU

    $txtSrc.Copy()
    $txt = $s7.Shapes.Paste().Item(1)
    $txt.Name = $set.txName
    $txt.Left = $set.txX / $EMU
    $txt.Top = $set.txY / $EMU
    $txt.Width = 790297 / $EMU
    $txt.Height = 276999 / $EMU
}

Write-Output "done"
